# Update "想去人数" (want-to-go count) values across sheets, matching the
# gh-pages generated-data refresh recorded in the commit history.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 9905
$wsExpo.Range("F5").Value = 588
$wsExpo.Range("F6").Value = 474

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 3

# Sheet "全部类型" (All types, combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9905
$wsAll.Range("F5").Value = 588
$wsAll.Range("F6").Value = 3
$wsAll.Range("F7").Value = 474
